$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells remain text (Excel would otherwise parse
# values like "4.00" or "35.015.65" as numbers/dates).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.015.65"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.850.07"
$ws.Range("E3").Value = "  +2.06%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.74"
$ws.Range("E5").Value = "  +3.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  +0.79%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.11"
$ws.Range("E8").Value = "  +6.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.328"
$ws.Range("E9").Value = "  +2.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0693"
$ws.Range("E10").Value = "  +1.89%  "

$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.120.89"
$ws.Range("E12").Value = "  +2.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.41"
$ws.Range("E13").Value = "  +1.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.844.55"
$ws.Range("E14").Value = "  +2.08%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.675"
$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.76"
$ws.Range("E16").Value = "  +4.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.050.01"
$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.17"
$ws.Range("E18").Value = "  +1.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0793"
$ws.Range("E19").Value = "  +1.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.01"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.16"
$ws.Range("E21").Value = "  +2.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.75"
$ws.Range("E22").Value = "  +2.64%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("E24").Value = "  +1.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.17"
$ws.Range("E25").Value = "  -1.67%  "

$ws.Range("E26").Value = "  +22.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.96"
$ws.Range("E27").Value = "  +3.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.62"
$ws.Range("E28").Value = "  +1.99%  "

$ws.Range("E29").Value = "  -0.82%  "

$ws.Range("E30").Value = "  +0.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0556"
$ws.Range("E31").Value = "  +2.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.00"
$ws.Range("E32").Value = "  +0.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.02"
$ws.Range("E33").Value = "  +3.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.66"
$ws.Range("E34").Value = "  +23.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.99"
$ws.Range("E35").Value = "  +10.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.31"
$ws.Range("E36").Value = "  +6.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.780"
$ws.Range("E37").Value = "  +13.77%  "

$ws.Range("E38").Value = "  +10.48%  "

$ws.Range("E39").Value = "  +5.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "90.04"
$ws.Range("E40").Value = "  -0.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.348.77"
$ws.Range("E41").Value = "  +1.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.69"
$ws.Range("E42").Value = "  +3.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.33"
$ws.Range("E43").Value = "  +2.57%  "

$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.42"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("B45").Value = "Gas"
$ws.Range("C45").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.44"
$ws.Range("E45").Value = "  +47.13%  "

$ws.Range("E46").Value = "  +7.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.74"
$ws.Range("E47").Value = "  -0.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.54"
$ws.Range("E48").Value = "  +6.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.038.36"
$ws.Range("E49").Value = "  +2.31%  "

$ws.Range("E50").Value = "  +2.72%  "

$ws.Range("E51").Value = "  +0.03%  "
